{"js": "// Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) in specific resume bullet\n// lines, matching the author's \"Implement quantitative metrics\n// highlighting across all resume formats\" commit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Map each target paragraph (matched by its exact, unique current text)\n// to the ordered list of metric substrings that must become bold + colored.\nconst targets = [\n  {\n    text:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    text:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n    metrics: [\"87%\", \"71%\", \"\\u00b14.2%\", \"\\u00b12.1%\"],\n  },\n  {\n    text:\n      \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"],\n  },\n  {\n    text:\n      \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"],\n  },\n  {\n    text:\n      \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"],\n  },\n];\n\n// Track how many paragraphs with a given text we've already matched, so\n// that duplicate bullet text (e.g. the two \"Achieved 87% ... 71%\" lines)\n// are each bound to their own distinct paragraph occurrence in document\n// order.\nconst usedCounts = new Map();\n\nfor (const target of targets) {\n  const seenSoFar = usedCounts.get(target.text) || 0;\n  let matchIndex = -1;\n  let seen = 0;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === target.text) {\n      if (seen === seenSoFar) {\n        matchIndex = i;\n        break;\n      }\n      seen++;\n    }\n  }\n  if (matchIndex === -1) {\n    throw new Error(\"Could not find target paragraph: \" + target.text);\n  }\n  usedCounts.set(target.text, seenSoFar + 1);\n\n  const paragraph = paragraphs.items[matchIndex];\n  for (const metric of target.metrics) {\n    const found = paragraph.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length === 0) {\n      throw new Error('Metric \"' + metric + '\" not found in paragraph: ' + target.text);\n    }\n    const range = found.items[0];\n    range.font.bold = true;\n    range.font.color = \"#2C3E50\";\n  }\n  await context.sync();\n}\n", "ps1": "# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) in specific resume bullet\n# lines, matching the author's \"Implement quantitative metrics\n# highlighting across all resume formats\" commit.\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n\n$targets = @(\n    @{\n        text    = \"$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        text    = \"$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + [char]0x00B1 + \"4.2% to \" + [char]0x00B1 + \"2.1%\"\n        metrics = @(\"87%\", \"71%\", ([char]0x00B1 + \"4.2%\"), ([char]0x00B1 + \"2.1%\"))\n    },\n    @{\n        text    = \"$bullet Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        metrics = @(\"1,200\")\n    },\n    @{\n        text    = \"$bullet Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        metrics = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        text    = \"$bullet Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        metrics = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        text    = \"$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        metrics = @(\"87%\", \"71%\")\n    }\n)\n\n# Track how many times each distinct paragraph text has already been\n# matched, so duplicate bullet text (the two \"Achieved 87% ... 71%\"\n# lines) are each bound to their own distinct paragraph occurrence, in\n# document order.\n$usedCounts = @{}\n\nforeach ($target in $targets) {\n    $wanted = 0\n    if ($usedCounts.ContainsKey($target.text)) {\n        $wanted = $usedCounts[$target.text]\n    }\n\n    $seen = 0\n    $matched = $null\n    foreach ($p in $d.Paragraphs) {\n        $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($ptext -eq $target.text) {\n            if ($seen -eq $wanted) {\n                $matched = $p\n                break\n            }\n            $seen = $seen + 1\n        }\n    }\n\n    if ($null -eq $matched) {\n        throw \"Could not find target paragraph: $($target.text)\"\n    }\n    $usedCounts[$target.text] = $wanted + 1\n\n    foreach ($m in $target.metrics) {\n        $r = $matched.Range\n        $find = $r.Find\n        $find.ClearFormatting()\n        $find.Text = $m\n        $find.MatchCase = $true\n        $found = $find.Execute()\n        if (-not $found) {\n            throw \"Metric '$m' not found in paragraph: $($target.text)\"\n        }\n        $r.Font.Bold = $true\n        $r.Font.Color = \"#2C3E50\"\n    }\n}\n"}
